# Generated by Katalon AI
# The "AI Generated" sheet previously held 4 locator columns
# (button_apiKeyAction_trNthChild / _1 / _2 / input_KeyName). Katalon AI
# regenerated the test data object so only the "input_KeyName" locator
# remains: columns A-C (the old trNthChild locators) are removed and the
# former column D (input_KeyName) slides left to become the new column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:C").Delete()
